# Actualiza la serie de intercambio comercial (M nov 2021) en la hoja de
# estacionalidad. Esto modifica:
#  - Las columnas de "media" (v_expo_media=W, v_impo_media=Y,
#    v_ic_media=AA) de los grupos de Noviembre y Diciembre, cuyo valor
#    depende del promedio historico y por tanto cambia para todas las
#    filas del grupo al incorporarse el dato real de nov-2021.
#  - Los valores "en crudo" (v_expo=V, v_impo=X, Intercambio_comercial=Z,
#    vX_var=AB, vM_var=AC) de las filas correspondientes al ultimo dato
#    observado de cada mes (Noviembre 2021 y Diciembre 2021).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grupo Noviembre (filas 112-122) ---
$wExpoMediaNov = 5323.570049492728
$yImpoMediaNov = 5090.122174681818
$aaIcMediaNov  = 233.4478748109091

112..121 | ForEach-Object {
    $ws.Range("W$_").Value  = $wExpoMediaNov
    $ws.Range("Y$_").Value  = $yImpoMediaNov
    $ws.Range("AA$_").Value = $aaIcMediaNov
}

# Fila 122 = Noviembre 2021: se actualiza el dato real observado.
$ws.Range("V122").Value  = 6191.188181310001
$ws.Range("W122").Value  = $wExpoMediaNov
$ws.Range("X122").Value  = 5767.0360416
$ws.Range("Y122").Value  = $yImpoMediaNov
$ws.Range("Z122").Value  = 424.15213971
$ws.Range("AA122").Value = $aaIcMediaNov
$ws.Range("AB122").Value = 0.3756873293367513
$ws.Range("AC122").Value = 0.4015497709303537

# --- Grupo Diciembre (filas 123-133) ---
$wExpoMediaDic = 5051.888272472727
$yImpoMediaDic = 4788.565382631818
$aaIcMediaDic  = 263.3228898409092

123..132 | ForEach-Object {
    $ws.Range("W$_").Value  = $wExpoMediaDic
    $ws.Range("Y$_").Value  = $yImpoMediaDic
    $ws.Range("AA$_").Value = $aaIcMediaDic
}

# Fila 133 = Diciembre 2021: se actualiza el dato real observado.
$ws.Range("V133").Value  = 6587.00000611
$ws.Range("W133").Value  = $wExpoMediaDic
$ws.Range("X133").Value  = 6215.73500319
$ws.Range("Y133").Value  = $yImpoMediaDic
$ws.Range("Z133").Value  = 371.26500292
$ws.Range("AA133").Value = $aaIcMediaDic
$ws.Range("AB133").Value = 0.8585985994147831
$ws.Range("AC133").Value = 0.5905395464668235
